# Add a "Price" column (C) with example product prices, as per
# commit: "add contoh product price in example"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("C1").Value = "Price"

# Example price for every product row (rows 2-26)
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 100000
}

# Mirror the author's final selection state as closely as the
# object model allows (multi-area selection with a non-primary
# active cell isn't reachable, so just land on B3).
[void]$ws.Range("B3").Select()
